$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (10:04 -> 11:04)
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 11:04"

# Update country rows: re-sorted order (name swaps) + refreshed case counts
# Row 8: Rusia
$ws.Cells.Item(8, 1).Value = "Rusia"
$ws.Cells.Item(8, 2).Value = 187859
$ws.Cells.Item(8, 3).Value = 10699
$ws.Cells.Item(8, 4).Value = 26608
$ws.Cells.Item(8, 5).Value = 159528
$ws.Cells.Item(8, 6).Value = 2300
$ws.Cells.Item(8, 7).Value = 98
$ws.Cells.Item(8, 8).Value = 1723

# Row 17: India
$ws.Cells.Item(17, 1).Value = "India"
$ws.Cells.Item(17, 2).Value = 56516
$ws.Cells.Item(17, 3).Value = 165
$ws.Cells.Item(17, 4).Value = 16867
$ws.Cells.Item(17, 5).Value = 37754
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 6
$ws.Cells.Item(17, 8).Value = 1895

# Row 29: Singapur
$ws.Cells.Item(29, 1).Value = "Singapur"
$ws.Cells.Item(29, 2).Value = 21707
$ws.Cells.Item(29, 3).Value = 768
$ws.Cells.Item(29, 4).Value = 1712
$ws.Cells.Item(29, 5).Value = 19975
$ws.Cells.Item(29, 6).Value = 19
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(29, 8).Value = 20

# Row 36: Polonia
$ws.Cells.Item(36, 1).Value = "Polonia"
$ws.Cells.Item(36, 2).Value = 15200
$ws.Cells.Item(36, 3).Value = 153
$ws.Cells.Item(36, 4).Value = 5184
$ws.Cells.Item(36, 5).Value = 9260
$ws.Cells.Item(36, 6).Value = 160
$ws.Cells.Item(36, 7).Value = 1
$ws.Cells.Item(36, 8).Value = 756

# Row 39: Banglades
$ws.Cells.Item(39, 1).Value = "Banglades"
$ws.Cells.Item(39, 2).Value = 13134
$ws.Cells.Item(39, 3).Value = 709
$ws.Cells.Item(39, 4).Value = 1910
$ws.Cells.Item(39, 5).Value = 11018
$ws.Cells.Item(39, 6).Value = 1
$ws.Cells.Item(39, 7).Value = 7
$ws.Cells.Item(39, 8).Value = 206

# Row 40: Indonesia
$ws.Cells.Item(40, 1).Value = "Indonesia"
$ws.Cells.Item(40, 2).Value = 12776
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 2381
$ws.Cells.Item(40, 5).Value = 9465
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 930

# Row 42: Filipinas
$ws.Cells.Item(42, 1).Value = "Filipinas"
$ws.Cells.Item(42, 2).Value = 10463
$ws.Cells.Item(42, 3).Value = 120
$ws.Cells.Item(42, 4).Value = 1734
$ws.Cells.Item(42, 5).Value = 8033
$ws.Cells.Item(42, 6).Value = 31
$ws.Cells.Item(42, 7).Value = 11
$ws.Cells.Item(42, 8).Value = 696

# Row 52: Australia
$ws.Cells.Item(52, 1).Value = "Australia"
$ws.Cells.Item(52, 2).Value = 6914
$ws.Cells.Item(52, 3).Value = 18
$ws.Cells.Item(52, 4).Value = 6079
$ws.Cells.Item(52, 5).Value = 738
$ws.Cells.Item(52, 6).Value = 20
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 97

# Row 54: Malasia
$ws.Cells.Item(54, 1).Value = "Malasia"
$ws.Cells.Item(54, 2).Value = 6535
$ws.Cells.Item(54, 3).Value = 68
$ws.Cells.Item(54, 4).Value = 4864
$ws.Cells.Item(54, 5).Value = 1564
$ws.Cells.Item(54, 6).Value = 18
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 107

# Row 66: Oman
$ws.Cells.Item(66, 1).Value = "Oman"
$ws.Cells.Item(66, 2).Value = 3112
$ws.Cells.Item(66, 3).Value = 154
$ws.Cells.Item(66, 4).Value = 1025
$ws.Cells.Item(66, 5).Value = 2072
$ws.Cells.Item(66, 6).Value = 17
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 8).Value = 15

# Row 67: Ghana
$ws.Cells.Item(67, 1).Value = "Ghana"
$ws.Cells.Item(67, 2).Value = 3091
$ws.Cells.Item(67, 3).Value = 0
$ws.Cells.Item(67, 4).Value = 303
$ws.Cells.Item(67, 5).Value = 2770
$ws.Cells.Item(67, 6).Value = 4
$ws.Cells.Item(67, 7).Value = 0
$ws.Cells.Item(67, 8).Value = 18

# Row 68: Armenia
$ws.Cells.Item(68, 1).Value = "Armenia"
$ws.Cells.Item(68, 2).Value = 3029
$ws.Cells.Item(68, 3).Value = 145
$ws.Cells.Item(68, 4).Value = 1218
$ws.Cells.Item(68, 5).Value = 1768
$ws.Cells.Item(68, 6).Value = 10
$ws.Cells.Item(68, 7).Value = 1
$ws.Cells.Item(68, 8).Value = 43

# Row 69: Tailandia
$ws.Cells.Item(69, 1).Value = "Tailandia"
$ws.Cells.Item(69, 2).Value = 3000
$ws.Cells.Item(69, 3).Value = 8
$ws.Cells.Item(69, 4).Value = 2784
$ws.Cells.Item(69, 5).Value = 161
$ws.Cells.Item(69, 6).Value = 61
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 55

# Row 82: Estonia
$ws.Cells.Item(82, 1).Value = "Estonia"
$ws.Cells.Item(82, 2).Value = 1725
$ws.Cells.Item(82, 3).Value = 5
$ws.Cells.Item(82, 4).Value = 704
$ws.Cells.Item(82, 5).Value = 965
$ws.Cells.Item(82, 6).Value = 4
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = 56

# Row 103: Sri Lanka
$ws.Cells.Item(103, 1).Value = "Sri Lanka"
$ws.Cells.Item(103, 2).Value = 824
$ws.Cells.Item(103, 3).Value = 1
$ws.Cells.Item(103, 4).Value = 240
$ws.Cells.Item(103, 5).Value = 575
$ws.Cells.Item(103, 6).Value = 1
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 9

# Row 150: Brunei
$ws.Cells.Item(150, 1).Value = "Brunei"
$ws.Cells.Item(150, 2).Value = 141
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 132
$ws.Cells.Item(150, 5).Value = 8
$ws.Cells.Item(150, 6).Value = 2
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 1

# Row 205: Montserrat
$ws.Cells.Item(205, 1).Value = "Montserrat"
$ws.Cells.Item(205, 2).Value = 11
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 7
$ws.Cells.Item(205, 5).Value = 3
$ws.Cells.Item(205, 6).Value = 1
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 1

# Row 206: Seychelles
$ws.Cells.Item(206, 1).Value = "Seychelles"
$ws.Cells.Item(206, 2).Value = 11
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 8
$ws.Cells.Item(206, 5).Value = 3
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0
